$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column C time-slot values for rows 2, 3, 6 and 7 (rows 4 and 5 stay the same)
$ws.Range("C2").Value = "9:30-9:35"
$ws.Range("C3").Value = "9:35-9:40"
$ws.Range("C6").Value = "19:25-19:30"
$ws.Range("C7").Value = "19:30-19:35"

# Update the active cell selection to C11
$ws.Range("C11").Select()
